$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Remove column B entirely (it duplicated column A in the prior edit)
$ws.Columns("B").Delete()

# Restore the original A3/A4 values ("c"/"d") that the prior edit had
# overwritten with duplicates of "a"/"b"
$ws.Range("A3").Value = "c"
$ws.Range("A4").Value = "d"

# Update selection to match final state (A5 selected)
$ws.Range("A5").Select()
